# Actualización 11 de Mayo - Mañana
# Updates the Aprobados / Reprobados / Por_Apro / Por_Repro / Promedio /
# Blancos / Por_Blan statistics for rows 2, 3 and 5 on the "2o Parcial"
# and "3er Parcial" sheets.

$wb = $excel.ActiveWorkbook

# --- "2o Parcial" sheet: rows 2, 3 and 5 (columns E:K) ---------------------
$ws2 = $wb.Worksheets.Item("2o Parcial")

# Row 2
$ws2.Range("E2").Value = 25
$ws2.Range("F2").Value = 14
$ws2.Range("G2").Value = 64.09999999999999
$ws2.Range("H2").Value = 35.9
$ws2.Range("I2").Value = 7.2
$ws2.Range("J2").Value = 0
$ws2.Range("K2").Value = 0

# Row 3
$ws2.Range("E3").Value = 30
$ws2.Range("F3").Value = 7
$ws2.Range("G3").Value = 81.08
$ws2.Range("H3").Value = 18.92
$ws2.Range("I3").Value = 8
$ws2.Range("J3").Value = 0
$ws2.Range("K3").Value = 0

# Row 5
$ws2.Range("E5").Value = 18
$ws2.Range("F5").Value = 6
$ws2.Range("G5").Value = 75
$ws2.Range("H5").Value = 25
$ws2.Range("I5").Value = 7.9
$ws2.Range("J5").Value = 0
$ws2.Range("K5").Value = 0

# --- "3er Parcial" sheet: rows 2, 3 and 5 -----------------------------------
$ws3 = $wb.Worksheets.Item("3er Parcial")

# Row 2
$ws3.Range("I2").Value = 7.5

# Row 3
$ws3.Range("E3").Value = 30
$ws3.Range("F3").Value = 7
$ws3.Range("G3").Value = 81.08
$ws3.Range("H3").Value = 18.92
$ws3.Range("I3").Value = 8.4

# Row 5
$ws3.Range("E5").Value = 18
$ws3.Range("F5").Value = 6
$ws3.Range("G5").Value = 75
$ws3.Range("H5").Value = 25
$ws3.Range("I5").Value = 8
